$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 4; $row++) {
    $ws.Cells.Item($row, 2).Value = "sports_club_coed"
    $ws.Cells.Item($row, 3).Value = "Disc Golf"
}
